# "#4 slides for defense slightly fixed"
#
# Slide 13 ("Tasks") had a bullet reading:
#   "Formulating the formal functional and non-functional requirements
#    for the software tool (the DB-nets Renew plugin)"
# The author trimmed it down to:
#   "Formulating the formal functional requirements for the software
#    tool (the DB-nets Renew plugin)"
# which ends up represented as three runs in the paragraph.

$p  = $ppt.ActivePresentation
$s  = $p.Slides.Item(13)
$sh = $s.Shapes.Item(5)
$tr = $sh.TextFrame.TextRange

# The bullet is the 5th paragraph of this text box.
$para = $tr.Paragraphs(5, 1)

# Remove "functional and non-" (chars 24..42, 1-based/length) from:
#   "Formulating the formal functional and non-functional requirements for the software tool (the DB-nets Renew plugin)"
# leaving:
#   "Formulating the formal functional requirements for the software tool (the DB-nets Renew plugin)"
$toDelete = $para.Characters(24, 19)
$toDelete.Text = ""

# Re-touch "functional requirements " (chars 24..47) in place so the
# paragraph ends up split into the same three runs the authored slide
# has ("Formulating the formal " / "functional requirements " / "for
# the software tool (the DB-nets Renew plugin)") while keeping the
# original run formatting.
$mid = $para.Characters(24, 24)
$mid.Text = $mid.Text
